$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph (paragraph 2 of the doc).
#    It sits right after the Heading1 title paragraph.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Eggomatic Free - Unique Slot Game
#    with Bonus Features") right before the final (DALLE prompt) paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r/>' +
              '<w:r><w:rPr><w:b/></w:rPr>' +
              '<w:t>Play Eggomatic Free - Unique Slot Game with Bonus Features</w:t></w:r>' +
              '</w:p>'
$newPara.Range.InsertXML($newParaXml)

# 3. Replace the DALLE image-prompt text (now the last paragraph) with the
#    meta-description sentence, keeping its existing italic run formatting.
$oldText = "Prompt: DALLE, please create a cartoon-style feature image for Eggomatic that incorporates a happy Maya warrior wearing glasses. The image should reflect the steampunk world of the game with pipes and futuristic machines in the background. The Maya warrior should be holding an egg while standing confidently in front of the EggOMatic machine. The overall tone of the image should be cheerful and inviting to entice players to try out the game."
$newText = "Read our Eggomatic slot game review. Play it free online with unique bonus features and stunning steampunk visuals."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
